$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Paragraph "The Greater London ... advice in accordance ... Charter."
#    - remove the inline _GoBack bookmark that currently sits between
#      "advice in" and " accordance ..."
#    - split the trailing text into its own run ("advice in" | " accordance...")
#    - insert a new, completely empty paragraph straight after it.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("advice in", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'advice in' paragraph"
}
$paraCharter = $r.Paragraphs(1)
$rngCharter = $paraCharter.Range

$xmlCharter = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:tabs><w:tab w:val="left" w:pos="5245"/></w:tabs>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr>
    <w:t xml:space="preserve">The Greater London Archaeological Advisory Service (GLAAS) provides archaeological </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr>
    <w:t>advice in</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr>
    <w:t xml:space="preserve"> accordance with the National Planning Policy Framework and GLAAS Charter.</w:t>
  </w:r>
</w:p>
"@
$rngCharter.InsertXML($xmlCharter)

# Insert a brand-new empty paragraph right after the paragraph we just rewrote.
$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Execute("GLAAS Charter.", $true, $false, $false, $false, $false, $true, 1, $false, "GLAAS Charter.^p", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "<Condition Type>" paragraph -> becomes "<Condition>" and gains the
#    _GoBack bookmark at its start; the ">" run now uses the Arial (cs) font
#    like the "<" -> "Condition" runs.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.ClearFormatting()
$found3 = $r3.Find.Execute("<Condition Type>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not find '<Condition Type>' paragraph"
}
$paraCond = $r3.Paragraphs(1)
$rngCond = $paraCond.Range
$xmlCond = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:tabs><w:tab w:val="left" w:pos="5245"/></w:tabs>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr>
    <w:t>&lt;</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:rFonts w:cs="Arial"/></w:rPr>
    <w:t>Condition</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:rFonts w:cs="Arial"/></w:rPr>
    <w:t>&gt;</w:t>
  </w:r>
</w:p>
"@
$rngCond.InsertXML($xmlCond)

# ---------------------------------------------------------------------------
# 3. The old, second "<Condition>" paragraph loses its text entirely
#    (becomes a blank paragraph, same as the one that follows it).
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.ClearFormatting()
$found4 = $r4.Find.Execute("<Condition>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found4) {
    throw "Could not find '<Condition>' paragraph"
}
$paraCond2 = $r4.Paragraphs(1)
$rngCond2 = $paraCond2.Range
$xmlCond2 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:tabs><w:tab w:val="left" w:pos="5245"/></w:tabs>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr>
  </w:pPr>
</w:p>
"@
$rngCond2.InsertXML($xmlCond2)

# ---------------------------------------------------------------------------
# 4. Remove the now-superfluous blank paragraph that used to trail the
#    second "<Condition>" paragraph (net paragraph count -1 for this block).
# ---------------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.ClearFormatting()
$r5.Find.Execute("<Condition>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraCond2b = $r5.Paragraphs(1)
$nextPara = $paraCond2b.Next()
$nextPara.Range.Delete()

# ---------------------------------------------------------------------------
# 5. Add <w:lastRenderedPageBreak/> immediately before the run
#    "This response relates solely to archaeological issues."
# ---------------------------------------------------------------------------
$r6 = $d.Content
$r6.Find.ClearFormatting()
$found6 = $r6.Find.Execute("This response relates solely to archaeological issues.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found6) {
    throw "Could not find the 'This response relates solely' paragraph"
}
$paraResp = $r6.Paragraphs(1)
$rngResp = $paraResp.Range
$xmlResp = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:tabs><w:tab w:val="left" w:pos="5245"/></w:tabs>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>This response relates solely to archaeological issues.</w:t>
  </w:r>
</w:p>
"@
$rngResp.InsertXML($xmlResp)
